$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old I11 cell (row 11 disappears entirely from the sheet)
$ws.Range("I11").ClearContents()

# New cells, set in the order that matches the shared-string table
# produced by the original authoring tool (index 16, 17, 18):
$ws.Range("I10").Value = "ответ:         1"
$ws.Range("D10").Value = "f (x,y,z)=a123xyz+a12xy+a13xz+a23YZ+a1x+a2y+a3z+a0"
$ws.Range("I1").Value = "X->(Y->Z)V((X->Y)vZ)=1"

# Column widths for H and I (values chosen so the engine's internal
# character-width rounding lands on the widths recorded in the workbook:
# H -> 11.5703125, I -> 36)
$ws.Columns.Item(8).ColumnWidth = 10.6666667
$ws.Columns.Item(9).ColumnWidth = 35.1666667
